$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, taken from the crypto-price refresh diff.
# NumberFormat is forced to Text ('@') before the write so numeric-looking
# strings (e.g. "4.615", "0.06170", "26.017.95") are stored verbatim as
# text instead of being parsed into floating-point numbers, then the style
# is reset to 'Normal' so no extra formatting is left behind on the cell.
$updates = [ordered]@{
    'D2' = '26.017.95'
    'E2' = '  +0.72%  '
    'D3' = '1.739.45'
    'E3' = '  -0.03%  '
    'E4' = '  +0.11%  '
    'D5' = '240.10'
    'E5' = '  +3.33%  '
    'E6' = '  +0.14%  '
    'D7' = '0.5289'
    'E7' = '  +2.28%  '
    'D8' = '0.2759'
    'E8' = '  -1.82%  '
    'D9' = '0.06170'
    'D10' = '1.740.83'
    'E10' = '  -0.85%  '
    'D11' = '0.07191'
    'E11' = '  +2.25%  '
    'D12' = '15.15'
    'E12' = '  -1.28%  '
    'D13' = '0.6428'
    'E13' = '  -0.96%  '
    'D14' = '4.615'
    'E14' = '  +1.90%  '
    'D15' = '77.61'
    'E15' = '  +0.71%  '
    'E16' = '  +0.15%  '
    'E17' = '  +0.13%  '
    'D18' = '26.033.27'
    'E18' = '  +0.83%  '
    'D19' = '11.79'
    'E19' = '  +2.71%  '
    'D20' = '0.000006782'
    'E20' = '  +2.66%  '
    'D21' = '1.964.73'
    'E21' = '  -0.59%  '
    'D22' = '4.372'
    'E22' = '  +5.70%  '
    'D23' = '8.640'
    'E23' = '  -0.03%  '
    'D24' = '5.262'
    'E24' = '  +2.18%  '
    'D25' = '140.31'
    'E25' = '  +0.75%  '
    'D26' = '1.515'
    'E26' = '  +0.05%  '
    'D27' = '15.26'
    'E27' = '  +1.21%  '
    'E28' = '  -2.42%  '
    'E29' = '  +3.65%  '
    'D30' = '0.08412'
    'E30' = '  +1.30%  '
    'D31' = '3.836'
    'E31' = '  +4.20%  '
    'D32' = '3.642'
    'E32' = '  +6.13%  '
    'D33' = '0.04595'
    'E33' = '  +2.11%  '
    'D34' = '2.653'
    'E34' = '  +1.60%  '
    'D35' = '0.9936'
    'E35' = '  +0.58%  '
    'D36' = '0.6246'
    'E36' = '  +1.36%  '
    'D37' = '2.700'
    'E37' = '  +1.72%  '
    'D38' = '0.01606'
    'E38' = '  +1.42%  '
    'D39' = '1.932'
    'E39' = '  -0.46%  '
    'E40' = '  +0.19%  '
    'D41' = '98.92'
    'E41' = '  -1.84%  '
    'D42' = '0.3888'
    'E42' = '  +1.00%  '
    'D43' = '0.7506'
    'E43' = '  +3.29%  '
    'D44' = '4.947'
    'E44' = '  -0.54%  '
    'E45' = '  +1.55%  '
    'D46' = '0.05311'
    'E46' = '  -1.81%  '
    'D47' = '6.228'
    'E47' = '  -1.05%  '
    'D48' = '54.74'
    'E48' = '  +2.95%  '
    'D49' = '30.80'
    'E49' = '  +2.96%  '
    'D50' = '0.3453'
    'E50' = '  +1.44%  '
    'D51' = '7.571'
    'E51' = '  -1.49%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$ref]
    $cell.Style = 'Normal'
}
